# "adding all DBMS work"
#
# The author briefly added a new slide (it was created right after the
# first slide, i.e. as the new slide #2) and then removed it again in the
# same editing session before the final save. The round trip is visible in
# the saved package only as bookkeeping (the slide id counter had already
# moved on to 484, relationship ids were bumped, coauthoring/"changesInfo"
# bookkeeping records an addSld followed by a delSld for sldId 484, etc.)
# -- the slide list itself ends up back at its original 100 slides, in
# their original order.
#
# Reproduce the same user action here: insert a new slide right after the
# first one, then delete it again.

$p = $ppt.ActivePresentation

$newSlide = $p.Slides.Add(2, 1)   # ppLayoutTitle; content is irrelevant, it gets removed below
$newSlide.Delete()
